# Weekly update: insert a new observation row for
# Hortaliza, Vega Monumental Concepción - Papa
# The new row is inserted before the current row 99, which pushes the
# existing rows 99-108 down to 100-109 (dimension grows from R108 to R109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 99 (shifts rows 99:108 -> 100:109)
$row = $ws.Rows.Item(99)
$row.Insert()

# Populate the newly inserted row 99 with the new record's data
$ws.Cells.Item(99, 1).Value  = 11                                   # A Mercado ID
$ws.Cells.Item(99, 2).Value  = "Vega Monumental Concepción"         # B Mercado
$ws.Cells.Item(99, 3).Value  = "Bíobío"                              # C Región
$ws.Cells.Item(99, 4).Value  = 44449                                 # D Fecha
$ws.Cells.Item(99, 5).Value  = 8                                     # E Codreg
$ws.Cells.Item(99, 6).Value  = 100114001                             # F Categoría ID
$ws.Cells.Item(99, 7).Value  = "Papa"                                # G Categoría
$ws.Cells.Item(99, 8).Value  = "Patagonia"                           # H Variedad
$ws.Cells.Item(99, 9).Value  = "1a (guarda)"                         # I Calidad
$ws.Cells.Item(99, 10).Value = 100                                   # J Volumen
$ws.Cells.Item(99, 11).Value = 8000                                  # K Precio mínimo
$ws.Cells.Item(99, 12).Value = 8500                                  # L Precio máximo
$ws.Cells.Item(99, 13).Value = 8250                                  # M Precio promedio ponderado
$ws.Cells.Item(99, 14).Value = "$/saco 25 kilos"                     # N Unidad de comercialización
$ws.Cells.Item(99, 15).Value = "Provincia de Arauco"                 # O Origen
$ws.Cells.Item(99, 16).Value = 330                                   # P Precio $/Kg
$ws.Cells.Item(99, 17).Value = 25                                    # Q Kg o Unidades
$ws.Cells.Item(99, 18).Value = "Hortaliza"                           # R Clasificación
